$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# New "Source" row (row 9)
$ws.Range("A9").Value = "Source"
$ws.Range("B9").Value = "ABS National Aboriginal and Torres Strait Islander Social Survey, various years and ABS Australian Aboriginal and Torres Strait Islander Health Survey, various years."

# New "References" row (row 10)
$ws.Range("A10").Value = "References"
$ws.Range("B10").Value = "Australian Housing and Urban Research Institute (AHURI), 2015, Indigenous housing - AHURI. [online] Available at: http://www.ahuri.edu.au/themes/indigenous_housing [Accessed 10 June 2015]."

# Formatting for the new body text cells - larger font, wrapped text (matches the
# other descriptive text cells in the sheet)
$ws.Range("B9:B10").Font.Size = 12
$ws.Range("B9:B10").WrapText = $true

# Row height adjustments
$ws.Rows(7).RowHeight = 20.95
$ws.Rows(8).RowHeight = 20.95
$ws.Rows(9).RowHeight = 26.95
$ws.Rows(10).RowHeight = 26.95

# Restore the selection to the updated range, keeping "Description" the active sheet
[void]$ws.Range("B7:B10").Select()
